# Add a header row ("Task:" / "Hours:") above the existing data, and append
# two new rows ("Pedestrians" and "Objective") at the bottom of the table.
# Also widen column B to fit the new "Hours:" header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all existing rows down by one to make room for the new header row.
$ws.Rows.Item(1).Insert()

# New rows appended at the end of the (now shifted) table.
$ws.Range("A19").Value = "Pedestrians"
$ws.Range("B19").Value = 0.5
$ws.Range("A20").Value = "Objective"
$ws.Range("B20").Value = 0.5

# New header row at the top.
$ws.Range("A1").Value = "Task:"
$ws.Range("B1").Value = "Hours:"

# Widen column B so the new "Hours:" header / longer values fit.
$ws.Columns.Item(2).ColumnWidth = 33.5

# Move the active selection the way the author left it.
$null = $ws.Range("A27").Select()
